# Generate Report for Handback
# The localization run picked up the "3f9aa1e4..." file's handback (it is now
# "Handed back: in sync with en-US" instead of "In Translation"), so the
# report rows get regenerated: rows are re-sorted by source file name and
# the handoff/handback bookkeeping columns for the three files rotate
# accordingly.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview" (table "Overview", A1:G4)
# Row order becomes: 3f9aa1e4..., ffff9b9a05c8..., fffffff18f5bde...
# ---------------------------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")

$ovw.Range("A2").Value = "3f9aa1e4-2994-4508-aca1-3fb142314ff7.md"
$ovw.Range("B2").Value = "e2e\3f9aa1e4-2994-4508-aca1-3fb142314ff7.md"

$ovw.Range("A3").Value = "ffff9b9a05c8-316b-403c-a34e-92095ea42bef.md"
$ovw.Range("B3").Value = "e2e\ffff9b9a05c8-316b-403c-a34e-92095ea42bef.md"

$ovw.Range("A4").Value = "fffffff18f5bde-2d4f-41d1-841b-b460142b2dea.md"
$ovw.Range("B4").Value = "e2e\fffffff18f5bde-2d4f-41d1-841b-b460142b2dea.md"
$ovw.Range("E4").Value = "Handed back: in sync with en-US"
$ovw.Range("F4").Value = "Handed back: in sync with en-US"
$ovw.Range("G4").Value = "2016-10-24 09:45:07"

foreach ($h in $ovw.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$B$2') { $h.TextToDisplay = "e2e\3f9aa1e4-2994-4508-aca1-3fb142314ff7.md" }
    elseif ($addr -eq '$B$3') { $h.TextToDisplay = "e2e\ffff9b9a05c8-316b-403c-a34e-92095ea42bef.md" }
    elseif ($addr -eq '$B$4') { $h.TextToDisplay = "e2e\fffffff18f5bde-2d4f-41d1-841b-b460142b2dea.md" }
}

# ---------------------------------------------------------------------
# Sheet "zh-cn" (table "zh_cn", A1:P4)
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "3f9aa1e4-2994-4508-aca1-3fb142314ff7.md"
$zh.Range("G2").Value = "3f9aa1e4-2994-4508-aca1-3fb142314ff7.a5a2129db27f9c907f997a41fa17f8f27e8c3ebe.zh-cn.xlf"
$zh.Range("H2").Value = "2016-10-24 09:50:41"
$zh.Range("I2").Value = "3f9aa1e4-2994-4508-aca1-3fb142314ff7.md"
$zh.Range("J2").Value = "3f9aa1e4-2994-4508-aca1-3fb142314ff7.a5a2129db27f9c907f997a41fa17f8f27e8c3ebe.zh-cn.xlf"
$zh.Range("K2").Value = "2016-10-24 09:51:28"

$zh.Range("A3").Value = "ffff9b9a05c8-316b-403c-a34e-92095ea42bef.md"
$zh.Range("F3").Value = "False"

$zh.Range("A4").Value = "fffffff18f5bde-2d4f-41d1-841b-b460142b2dea.md"
$zh.Range("C4").Value = "Handed back: in sync with en-US"
$zh.Range("F4").Value = "True"
$zh.Range("G4").Value = "5e1a45c3-46ab-4f7b-83e1-c01e94c7b632.868128d0013e40b92172d4f8000cc902d6794320.zh-cn.xlf"
$zh.Range("H4").Value = "2016-10-24 09:44:55"
$zh.Range("I4").Value = "5e1a45c3-46ab-4f7b-83e1-c01e94c7b632.md"
$zh.Range("J4").Value = "5e1a45c3-46ab-4f7b-83e1-c01e94c7b632.868128d0013e40b92172d4f8000cc902d6794320.zh-cn.xlf"
$zh.Range("K4").Value = "2016-10-24 09:45:36"
$zh.Range("P4").Value = ""

foreach ($h in $zh.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') { $h.TextToDisplay = "3f9aa1e4-2994-4508-aca1-3fb142314ff7.md" }
    elseif ($addr -eq '$I$2') { $h.TextToDisplay = "3f9aa1e4-2994-4508-aca1-3fb142314ff7.md" }
    elseif ($addr -eq '$A$3') { $h.TextToDisplay = "ffff9b9a05c8-316b-403c-a34e-92095ea42bef.md" }
    elseif ($addr -eq '$A$4') { $h.TextToDisplay = "fffffff18f5bde-2d4f-41d1-841b-b460142b2dea.md" }
    elseif ($addr -eq '$I$4') { $h.TextToDisplay = "5e1a45c3-46ab-4f7b-83e1-c01e94c7b632.md" }
}

# ---------------------------------------------------------------------
# Sheet "de-de" (table "de_de", A1:P4)
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "3f9aa1e4-2994-4508-aca1-3fb142314ff7.md"
$de.Range("G2").Value = "3f9aa1e4-2994-4508-aca1-3fb142314ff7.a5a2129db27f9c907f997a41fa17f8f27e8c3ebe.de-de.xlf"
$de.Range("H2").Value = "2016-10-24 09:50:52"
$de.Range("I2").Value = "3f9aa1e4-2994-4508-aca1-3fb142314ff7.md"
$de.Range("J2").Value = "3f9aa1e4-2994-4508-aca1-3fb142314ff7.a5a2129db27f9c907f997a41fa17f8f27e8c3ebe.de-de.xlf"
$de.Range("K2").Value = "2016-10-24 09:51:45"

$de.Range("A3").Value = "ffff9b9a05c8-316b-403c-a34e-92095ea42bef.md"
$de.Range("F3").Value = "False"

$de.Range("A4").Value = "fffffff18f5bde-2d4f-41d1-841b-b460142b2dea.md"
$de.Range("C4").Value = "Handed back: in sync with en-US"
$de.Range("F4").Value = "True"
$de.Range("G4").Value = "5e1a45c3-46ab-4f7b-83e1-c01e94c7b632.868128d0013e40b92172d4f8000cc902d6794320.de-de.xlf"
$de.Range("H4").Value = "2016-10-24 09:45:07"
$de.Range("I4").Value = "5e1a45c3-46ab-4f7b-83e1-c01e94c7b632.md"
$de.Range("J4").Value = "5e1a45c3-46ab-4f7b-83e1-c01e94c7b632.868128d0013e40b92172d4f8000cc902d6794320.de-de.xlf"
$de.Range("K4").Value = "2016-10-24 09:45:52"
$de.Range("P4").Value = ""

foreach ($h in $de.Hyperlinks) {
    $addr = $h.Range.Address()
    if ($addr -eq '$A$2') { $h.TextToDisplay = "3f9aa1e4-2994-4508-aca1-3fb142314ff7.md" }
    elseif ($addr -eq '$I$2') { $h.TextToDisplay = "3f9aa1e4-2994-4508-aca1-3fb142314ff7.md" }
    elseif ($addr -eq '$A$3') { $h.TextToDisplay = "ffff9b9a05c8-316b-403c-a34e-92095ea42bef.md" }
    elseif ($addr -eq '$A$4') { $h.TextToDisplay = "fffffff18f5bde-2d4f-41d1-841b-b460142b2dea.md" }
    elseif ($addr -eq '$I$4') { $h.TextToDisplay = "5e1a45c3-46ab-4f7b-83e1-c01e94c7b632.md" }
}

# ---------------------------------------------------------------------
# Column P on the zh-cn / de-de sheets shrinks now that the long error
# message is gone (Excel auto-fits the column after the edit).
# ---------------------------------------------------------------------
$zh.Columns.Item(16).ColumnWidth = 13.7470531463623
$de.Columns.Item(16).ColumnWidth = 13.7470531463623

Write-Output "done"
